$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.651.37"
$ws.Range("E2").Value = "  +0.89%  "

$ws.Range("D3").Value = "2.611.89"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'600.98"
$ws.Range("E5").Value = "  +1.13%  "

$ws.Range("D6").Value = "'154.24"
$ws.Range("E6").Value = "  -0.33%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'0.548"
$ws.Range("E8").Value = "  +1.37%  "

$ws.Range("D9").Value = "2.609.83"
$ws.Range("E9").Value = "  +0.26%  "

$ws.Range("E10").Value = "  +7.21%  "

$ws.Range("E12").Value = "  +0.36%  "

$ws.Range("D13").Value = "'0.353"
$ws.Range("E13").Value = "  -1.56%  "

$ws.Range("D14").Value = "'28.02"
$ws.Range("E14").Value = "  -0.57%  "

$ws.Range("E15").Value = "  +2.10%  "

$ws.Range("D16").Value = "3.086.58"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").Value = "67.553.14"
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("D18").Value = "2.609.34"
$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").Value = "'11.23"
$ws.Range("E19").Value = "  -1.79%  "

$ws.Range("D20").Value = "'364.99"
$ws.Range("E20").Value = "  +2.95%  "

$ws.Range("E21").Value = "  -3.34%  "

$ws.Range("E22").Value = "  -0.68%  "

$ws.Range("D23").Value = "'2.08"
$ws.Range("E23").Value = "  +1.61%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("E25").Value = "  +3.54%  "

$ws.Range("E26").Value = "  -4.32%  "

$ws.Range("D27").Value = "'0.0000104"
$ws.Range("E27").Value = "  +1.01%  "

$ws.Range("D28").Value = "2.751.29"
$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("D29").Value = "'579.49"
$ws.Range("E29").Value = "  -2.85%  "

$ws.Range("E30").Value = "  +0.29%  "

$ws.Range("D31").Value = "'1.42"
$ws.Range("E31").Value = "  -2.69%  "

$ws.Range("D32").Value = "'7.91"
$ws.Range("E32").Value = "  -2.35%  "

$ws.Range("E33").Value = "  -0.72%  "

$ws.Range("E34").Value = "  -2.57%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("E36").Value = "  -3.26%  "

$ws.Range("D37").Value = "'4.95"
$ws.Range("E37").Value = "  -1.70%  "

$ws.Range("D38").Value = "'19.37"
$ws.Range("E38").Value = "  +0.16%  "

$ws.Range("D39").Value = "'155.43"
$ws.Range("E39").Value = "  +0.96%  "

$ws.Range("D40").Value = "'0.371"
$ws.Range("E40").Value = "  +0.24%  "

$ws.Range("D41").Value = "'5.37"
$ws.Range("E41").Value = "  -1.79%  "

$ws.Range("E42").Value = "  +1.50%  "

$ws.Range("D43").Value = "'2.64"
$ws.Range("E43").Value = "  +0.39%  "

$ws.Range("D44").Value = "'41.10"
$ws.Range("E44").Value = "  -0.96%  "

$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("E46").Value = "  -0.27%  "

$ws.Range("D47").Value = "'156.06"
$ws.Range("E47").Value = "  -0.53%  "

$ws.Range("D48").Value = "0.0₆0284"
$ws.Range("E48").Value = "  -8.65%  "

$ws.Range("D49").Value = "'3.75"
$ws.Range("E49").Value = "  -0.84%  "

$ws.Range("D50").Value = "'20.92"
$ws.Range("E50").Value = "  -2.00%  "

$ws.Range("E51").Value = "  -0.18%  "
